$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.983.24"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.172.70"
$ws.Range("E3").Value = "  -2.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.76"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.03"
$ws.Range("E7").Value = "  -6.35%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.566"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.22"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0926"
$ws.Range("E11").Value = "  -4.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.51"
$ws.Range("E12").Value = "  -16.18%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "2.496.40"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.16"
$ws.Range("E17").Value = "  -6.14%  "
$ws.Range("D18").Value = "2.195.92"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").Value = "40.943.24"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.46"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.15"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  -8.80%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.36"
$ws.Range("E27").Value = "  +11.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.41"
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.89"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  -11.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.22"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.68"
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0737"
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "25.57"
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.07"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("E39").Value = "  +7.26%  "
$ws.Range("E40").Value = "  -6.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.49"
$ws.Range("E41").Value = "  -9.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "61.24"
$ws.Range("E42").Value = "  -12.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.38"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.78"
$ws.Range("E44").Value = "  -6.71%  "
$ws.Range("E45").Value = "  -11.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.52"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0983"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  -4.03%  "
$ws.Range("E51").Value = "  -0.43%  "
